$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column A, shifting B:F left to A:E
$ws.Columns.Item(1).Delete()

# Fix the renamed header string (MODEL_CONDITION -> MODELCONDITION), now in D1
$ws.Cells.Replace("MODEL_CONDITION", "MODELCONDITION") | Out-Null
